$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.871.12'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '3.816.17'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '626.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '3.814.34'
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.454'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '4.456.74'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '3.820.74'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '68.881.59'
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '465.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '3.968.01'
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("E31").Value = '  -3.70%  '
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("E38").Value = '  +7.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.99%  '
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.977'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  +2.85%  '
$ws.Range("E45").Value = '  +0.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '154.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.40%  '
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("E50").Value = '  +2.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '380.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.63%  '
